# Apply updates to column F (dSF) for specific rows as described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 1
    9  = 1
    14 = 0
    23 = 0
    26 = 4
    31 = 0
    33 = 1
    35 = 0
    38 = 0
    39 = 3
    40 = -2
    42 = 1
    45 = 0
    47 = -8
    49 = -2
    50 = 2
    58 = -1
    61 = 0
    68 = -2
    77 = 1
    84 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$wb.Save()
